$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# New header cell F1 - copy formatting from E1 (the existing header style),
# then overwrite its value with the new column name.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Data values for the new time_taken column (rows 2-8), unstyled like the
# rest of the data cells.
$times = @(
    "2021-10-05 13:39:12.256561",
    "2021-10-05 13:39:12.256573",
    "2021-10-05 13:39:12.256577",
    "2021-10-05 13:39:12.256580",
    "2021-10-05 13:39:12.256584",
    "2021-10-05 13:39:12.256587",
    "2021-10-05 13:39:12.256590"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
